# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45177 (2023-09-08) to 45178 (2023-09-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = 45177
$newValue = 45178

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 176
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
